$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (491) down to the
# new rows (492:498) so the new cells inherit the same cell styles used
# throughout the table (date format, data font, etc.).
$ws.Range("A491:I491").Copy() | Out-Null
$ws.Range("A492:I498").PasteSpecial(-4122) | Out-Null

# New training-log entries (date 2025-10-21 -> serial 45951)
$rows = @(
    @{ Row=492; Date=45951; Name="Yoann Martelat"; Volume=70; Intensite=3; Fatigue=4; Douleur=4; Loc="Genou";              Plaisir=8 },
    @{ Row=493; Date=45951; Name="Kamal Bafounta";  Volume=70; Intensite=6; Fatigue=3; Douleur=4; Loc="Genou";              Plaisir=10 },
    @{ Row=494; Date=45951; Name="Maé Clavel";      Volume=70; Intensite=4; Fatigue=4; Douleur=0; Loc=$null;                Plaisir=6 },
    @{ Row=495; Date=45951; Name="Levy Ndoutoume";  Volume=70; Intensite=7; Fatigue=7; Douleur=1; Loc="Ischio";             Plaisir=8 },
    @{ Row=496; Date=45951; Name="Naim Ighbane";    Volume=70; Intensite=4; Fatigue=4; Douleur=0; Loc=$null;                Plaisir=0 },
    @{ Row=497; Date=45951; Name="Amir Etien";      Volume=70; Intensite=5; Fatigue=5; Douleur=6; Loc="Ischio";             Plaisir=6 },
    @{ Row=498; Date=45951; Name="Naim Dhib";       Volume=70; Intensite=5; Fatigue=5; Douleur=3; Loc="Courbaturé hanche";  Plaisir=2 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Volume
    $ws.Cells.Item($row, 4).Value = $r.Intensite
    $ws.Cells.Item($row, 5).Value = $r.Fatigue
    $ws.Cells.Item($row, 6).Value = $r.Douleur
    if ($r.Loc) {
        $ws.Cells.Item($row, 7).Value = $r.Loc
    } else {
        # Rows with no pain location keep the "empty" look-and-feel used
        # elsewhere in the sheet (style taken from another blank cell, e.g. G486)
        $ws.Range("G486").Copy() | Out-Null
        $ws.Range("G$row").PasteSpecial(-4122) | Out-Null
        $ws.Cells.Item($row, 7).ClearContents() | Out-Null
    }
    $ws.Cells.Item($row, 8).Value = $r.Plaisir
}

# Charge = Volume * Intensite, same formula used throughout column I
$ws.Range("I492:I498").Formula = "=C492*D492"

# Restore the view scroll position / active selection as left by the author
$excel.ActiveWindow.ScrollRow = 479
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K494").Select() | Out-Null
